$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.973.64'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.48%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.643.69'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.44%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.88%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.27'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.25%  '
$ws.Range('E6').Value = '  -0.86%  '
$ws.Range('E7').Value = '  -0.84%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2574'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.71%  '
$ws.Range('E9').Value = '  -0.25%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.68'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.28%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07762'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.55%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.272'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.07%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.645.14'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.18%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.870.11'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5473'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.29%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₅7950'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.70%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.50'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.97%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '25.995.16'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.49%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.000'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '203.19'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.392'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.03%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.910'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.38%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.002'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.39%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.001'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.88%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.881'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.46%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '140.99'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.56%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1139'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.53%  '
$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.828'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.37%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.73'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.50%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.240'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.32%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.04937'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.17%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.274'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.18%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.219'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.53%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.545'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.35%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.360'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.43%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.8944'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.29%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.618'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.65%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.150.47'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.03%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5595'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.77%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01567'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.43%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.000'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.93%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.710'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.72%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8069'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.36%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.77'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.48%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.780.59'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.49%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0₈117'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.38%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4512'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.68%  '
$ws.Range('E48').Value = '  -0.66%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '54.75'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.77%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05043'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.69%  '
$ws.Range('E51').Value = '  -0.72%  '
